$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- NroCuenta (account number) updated for both data rows ---
$ws.Range("E2").Value = 2240451788
$ws.Range("E3").Value = 2240451788

# --- FechaInicio (row 2 plain text; row 3 forced text that keeps the leading apostrophe) ---
$ws.Range("N2").Value = "'07/04/2021"
$ws.Range("N3").Value = "''07/04/2021"

# --- Row 2: Patente / Motor / Chasis bumped to the new RGA007 plate family ---
$ws.Range("V2").Value = "RGA007"
$ws.Range("W2").Value = "ABC12RGA007"
$ws.Range("X2").Value = "ZAZ123RGA007"

# --- Row 3: Patente / Motor / Chasis bumped to the new RGA008 plate family ---
$ws.Range("V3").Value = "RGA008"
$ws.Range("W3").Value = "ABC12RGA008"
$ws.Range("X3").Value = "ZAZ123RGA008"

# --- restore the view's last selection/scroll state ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 15
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("Z8").Select()
